$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("D6").Value = 44434
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 28000
$ws.Range("L6").Value = 30000
$ws.Range("M6").Value = 29000
$ws.Range("P6").Value = 1160

# Row 7
$ws.Range("D7").Value = 44384
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 26000
$ws.Range("L7").Value = 28000
$ws.Range("M7").Value = 27000
$ws.Range("P7").Value = 1080

# Row 8
$ws.Range("D8").Value = 44363
$ws.Range("J8").Value = 240
$ws.Range("K8").Value = 28000
$ws.Range("L8").Value = 30000
$ws.Range("M8").Value = 29000
$ws.Range("P8").Value = 1160

# Row 9
$ws.Range("D9").Value = 44349
$ws.Range("J9").Value = 600

# Row 10
$ws.Range("D10").Value = 44385
$ws.Range("J10").Value = 500
$ws.Range("K10").Value = 26000
$ws.Range("L10").Value = 28000
$ws.Range("M10").Value = 27000
$ws.Range("P10").Value = 1080

# Row 11
$ws.Range("D11").Value = 44427
$ws.Range("J11").Value = 300
$ws.Range("K11").Value = 28000
$ws.Range("L11").Value = 30000
$ws.Range("M11").Value = 29000
$ws.Range("P11").Value = 1160

# Row 12
$ws.Range("D12").Value = 44413
$ws.Range("J12").Value = 700

# Row 13
$ws.Range("D13").Value = 44377
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 26000
$ws.Range("L13").Value = 28000
$ws.Range("M13").Value = 27000
$ws.Range("P13").Value = 1080

# Row 14
$ws.Range("D14").Value = 44426
$ws.Range("J14").Value = 400
$ws.Range("K14").Value = 28000
$ws.Range("L14").Value = 30000
$ws.Range("M14").Value = 29000
$ws.Range("P14").Value = 1160

# Row 15
$ws.Range("D15").Value = 44412
$ws.Range("J15").Value = 600
$ws.Range("K15").Value = 25000
$ws.Range("L15").Value = 27000
$ws.Range("M15").Value = 26000
$ws.Range("P15").Value = 1040

# Row 16
$ws.Range("D16").Value = 44371
$ws.Range("J16").Value = 500

# Row 17
$ws.Range("D17").Value = 44364
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 28000
$ws.Range("L17").Value = 30000
$ws.Range("M17").Value = 29000
$ws.Range("P17").Value = 1160

# Row 18
$ws.Range("D18").Value = 44435
$ws.Range("J18").Value = 900
$ws.Range("K18").Value = 28000
$ws.Range("L18").Value = 30000
$ws.Range("M18").Value = 29000
$ws.Range("P18").Value = 1160

# Row 19
$ws.Range("D19").Value = 44405
$ws.Range("J19").Value = 500
$ws.Range("K19").Value = 26000
$ws.Range("L19").Value = 28000
$ws.Range("M19").Value = 27000
$ws.Range("P19").Value = 1080

# Row 20
$ws.Range("D20").Value = 44391
$ws.Range("J20").Value = 100

# Row 21
$ws.Range("D21").Value = 44350
$ws.Range("J21").Value = 700
$ws.Range("K21").Value = 28000
$ws.Range("L21").Value = 30000
$ws.Range("M21").Value = 29000
$ws.Range("P21").Value = 1160

# Row 22
$ws.Range("D22").Value = 44399
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 26000
$ws.Range("L22").Value = 28000
$ws.Range("M22").Value = 27000
$ws.Range("P22").Value = 1080

# Row 23
$ws.Range("D23").Value = 44419
$ws.Range("J23").Value = 600
$ws.Range("L23").Value = 29000
$ws.Range("M23").Value = 28000
$ws.Range("P23").Value = 1120

# Row 24
$ws.Range("D24").Value = 44420
$ws.Range("J24").Value = 700
$ws.Range("K24").Value = 27000
$ws.Range("L24").Value = 29000
$ws.Range("M24").Value = 28000
$ws.Range("P24").Value = 1120

# Row 25
$ws.Range("D25").Value = 44433
$ws.Range("J25").Value = 400
$ws.Range("K25").Value = 28000
$ws.Range("L25").Value = 30000
$ws.Range("M25").Value = 29000
$ws.Range("P25").Value = 1160

# Row 26
$ws.Range("A26").Value = 2
$ws.Range("B26").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C26").Value = 'Coquimbo'
$ws.Range("D26").Value = 44370
$ws.Range("D26").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E26").Value = 4
$ws.Range("F26").Value = 100112022
$ws.Range("G26").Value = 'Arveja Verde'
$ws.Range("H26").Value = 'Perfection'
$ws.Range("I26").Value = 'Primera'
$ws.Range("J26").Value = 400
$ws.Range("K26").Value = 27000
$ws.Range("L26").Value = 28000
$ws.Range("M26").Value = 27500
$ws.Range("N26").Value = '$/malla 25 kilos'
$ws.Range("O26").Value = 'Provincia de Limarí'
$ws.Range("P26").Value = 1100
$ws.Range("Q26").Value = 25
$ws.Range("R26").Value = 'Hortaliza'

# Row 27
$ws.Range("A27").Value = 2
$ws.Range("B27").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C27").Value = 'Coquimbo'
$ws.Range("D27").Value = 44398
$ws.Range("D27").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E27").Value = 4
$ws.Range("F27").Value = 100112022
$ws.Range("G27").Value = 'Arveja Verde'
$ws.Range("H27").Value = 'Perfection'
$ws.Range("I27").Value = 'Primera'
$ws.Range("J27").Value = 500
$ws.Range("K27").Value = 26000
$ws.Range("L27").Value = 28000
$ws.Range("M27").Value = 27000
$ws.Range("N27").Value = '$/malla 25 kilos'
$ws.Range("O27").Value = 'Provincia de Limarí'
$ws.Range("P27").Value = 1080
$ws.Range("Q27").Value = 25
$ws.Range("R27").Value = 'Hortaliza'

# Row 28
$ws.Range("A28").Value = 2
$ws.Range("B28").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C28").Value = 'Coquimbo'
$ws.Range("D28").Value = 44343
$ws.Range("D28").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E28").Value = 4
$ws.Range("F28").Value = 100112022
$ws.Range("G28").Value = 'Arveja Verde'
$ws.Range("H28").Value = 'Perfection'
$ws.Range("I28").Value = 'Primera'
$ws.Range("J28").Value = 200
$ws.Range("K28").Value = 26000
$ws.Range("L28").Value = 28000
$ws.Range("M28").Value = 27000
$ws.Range("N28").Value = '$/malla 25 kilos'
$ws.Range("O28").Value = 'Provincia de Limarí'
$ws.Range("P28").Value = 1080
$ws.Range("Q28").Value = 25
$ws.Range("R28").Value = 'Hortaliza'
